$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '58.237.96'
$ws.Range('E2').Value = '  -1.18%  '
Set-TextValue $ws.Range('D3') '2.486.23'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue $ws.Range('D5') '521.00'
$ws.Range('E5').Value = '  -2.32%  '
Set-TextValue $ws.Range('D6') '134.63'
$ws.Range('E6').Value = '  -0.68%  '
Set-TextValue $ws.Range('D7') '0.996'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('E8').Value = '  -1.52%  '
Set-TextValue $ws.Range('D9') '2.504.02'
$ws.Range('E9').Value = '  -0.09%  '
Set-TextValue $ws.Range('D10') '0.0989'
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('E13').Value = '  -2.09%  '
Set-TextValue $ws.Range('D14') '2.928.32'
$ws.Range('E14').Value = '  -0.62%  '
Set-TextValue $ws.Range('D15') '58.169.30'
$ws.Range('E15').Value = '  -1.19%  '
Set-TextValue $ws.Range('D16') '22.14'
$ws.Range('E16').Value = '  -2.83%  '
$ws.Range('E17').Value = '  -1.75%  '
Set-TextValue $ws.Range('D18') '2.495.13'
$ws.Range('E18').Value = '  -0.08%  '
Set-TextValue $ws.Range('D19') '10.69'
$ws.Range('E19').Value = '  -3.20%  '
Set-TextValue $ws.Range('D20') '321.99'
$ws.Range('E20').Value = '  -0.41%  '
Set-TextValue $ws.Range('D21') '4.18'
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  -2.95%  '
Set-TextValue $ws.Range('D24') '64.76'
$ws.Range('E24').Value = '  -0.38%  '
Set-TextValue $ws.Range('D25') '0.412'
$ws.Range('E25').Value = '  -1.61%  '
Set-TextValue $ws.Range('D26') '0.162'
$ws.Range('E26').Value = '  -1.36%  '
Set-TextValue $ws.Range('D27') '0.992'
$ws.Range('E27').Value = '  -0.78%  '
Set-TextValue $ws.Range('D28') '7.41'
$ws.Range('E28').Value = '  -1.53%  '
Set-TextValue $ws.Range('D29') '0.0₃0751'
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('E30').Value = '  -0.26%  '
Set-TextValue $ws.Range('D31') '6.35'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('E33').Value = '  +5.02%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -0.21%  '
Set-TextValue $ws.Range('D36') '18.11'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('E38').Value = '  -0.25%  '
Set-TextValue $ws.Range('D39') '36.67'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('E40').Value = '  -2.65%  '
Set-TextValue $ws.Range('D41') '0.800'
$ws.Range('E41').Value = '  +0.13%  '
Set-TextValue $ws.Range('D42') '276.03'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('E43').Value = '  -3.31%  '
Set-TextValue $ws.Range('D44') '5.07'
$ws.Range('E44').Value = '  +1.93%  '
Set-TextValue $ws.Range('D45') '0.599'
$ws.Range('E45').Value = '  -0.55%  '
Set-TextValue $ws.Range('D46') '124.11'
$ws.Range('E46').Value = '  -3.54%  '
Set-TextValue $ws.Range('D47') '0.0910'
$ws.Range('E47').Value = '  -1.61%  '
Set-TextValue $ws.Range('D48') '0.0492'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('E49').Value = '  -1.61%  '
Set-TextValue $ws.Range('D50') '17.16'
$ws.Range('E50').Value = '  -0.43%  '
Set-TextValue $ws.Range('D51') '1.740.73'
$ws.Range('E51').Value = '  -0.50%  '
